# ============================================================================
# Kujata_Profits workbook update — scheduled runner refresh of market pricing
# and leve-profit computations (currentAveragePrice*, LevePrice*, LeveProfit*).
#
# This script writes the refreshed static values for each affected leve row on
# every job sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), matching a refresh of
# the underlying market-board data. A few rows also gain/lose a profit cell
# (e.g. when a value becomes newly computable, or no longer applicable).
# ============================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: ALC
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 2394.7778
$ws.Range("I58").Value = 844.1667
$ws.Range("J58").Value = 3170.0833
$ws.Range("K58").Value = 2532.5001
$ws.Range("L58").Value = 9510.249899999999
$ws.Range("M58").Value = -2382.5001
$ws.Range("N58").Value = -9810.249899999999
$ws.Range("H106").Value = 5536.9287
$ws.Range("I106").Value = 6081.36
$ws.Range("K106").Value = 6081.36
$ws.Range("M106").Value = -5450.36
$ws.Range("H138").Value = 2718.7625
$ws.Range("I138").Value = 2577.3333
$ws.Range("J138").Value = 2736.6902
$ws.Range("K138").Value = 7731.999899999999
$ws.Range("L138").Value = 8210.070599999999
$ws.Range("M138").Value = -2591.999899999999
$ws.Range("N138").Value = -18490.0706

# ---------------------------------------------------------------------------
# Sheet: ARM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 50
$ws.Range("I5").Value = 50
$ws.Range("K5").Value = 50
$ws.Range("M5").Value = 62
$ws.Range("H45").Value = 1333.1666
$ws.Range("I45").Value = 1203
$ws.Range("J45").Value = 1437.3
$ws.Range("K45").Value = 1203
$ws.Range("L45").Value = 1437.3
$ws.Range("M45").Value = -826
$ws.Range("N45").Value = -2191.3
$ws.Range("H61").Value = 58824896
$ws.Range("I61").Value = 83334104
$ws.Range("K61").Value = 83334104
$ws.Range("M61").Value = -83333892
$ws.Range("H63").Value = 2075.1333
$ws.Range("I63").Value = 1945.375
$ws.Range("J63").Value = 2594.1667
$ws.Range("K63").Value = 1945.375
$ws.Range("L63").Value = 2594.1667
$ws.Range("M63").Value = -1259.375
$ws.Range("N63").Value = -3966.1667
$ws.Range("H66").Value = 2075.1333
$ws.Range("I66").Value = 1945.375
$ws.Range("J66").Value = 2594.1667
$ws.Range("K66").Value = 9726.875
$ws.Range("L66").Value = 12970.8335
$ws.Range("M66").Value = -6294.875
$ws.Range("N66").Value = -19834.8335
$ws.Range("H74").Value = 2630
$ws.Range("I74").Value = 1516.25
$ws.Range("K74").Value = 1516.25
$ws.Range("M74").Value = -642.25
$ws.Range("H77").Value = 2630
$ws.Range("I77").Value = 1516.25
$ws.Range("K77").Value = 7581.25
$ws.Range("M77").Value = -3213.25
$ws.Range("H81").Value = 46666.668
$ws.Range("J81").Value = 46666.668
$ws.Range("L81").Value = 46666.668
$ws.Range("N81").Value = -48662.668
$ws.Range("H84").Value = 46666.668
$ws.Range("J84").Value = 46666.668
$ws.Range("L84").Value = 140000.004
$ws.Range("N84").Value = -149984.004
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H88").Value = 1833.3334
$ws.Range("I88").Value = 2000
$ws.Range("J88").Value = 1750
$ws.Range("K88").Value = 2000
$ws.Range("L88").Value = 1750
$ws.Range("M88").Value = -1594
$ws.Range("N88").Value = -2562
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H91").Value = 1833.3334
$ws.Range("I91").Value = 2000
$ws.Range("J91").Value = 1750
$ws.Range("K91").Value = 2000
$ws.Range("L91").Value = 1750
$ws.Range("M91").Value = -596
$ws.Range("N91").Value = -4558
$ws.Range("H93").Value = 31000
$ws.Range("J93").Value = 31000
$ws.Range("L93").Value = 31000
$ws.Range("N93").Value = -35992
$ws.Range("H94").Value = 16000
$ws.Range("J94").Value = 16000
$ws.Range("L94").Value = 16000
$ws.Range("N94").Value = -17802
$ws.Range("H132").Value = 2981.1538
$ws.Range("I132").Value = 2594.8215
$ws.Range("J132").Value = 3964.5454
$ws.Range("K132").Value = 7784.4645
$ws.Range("L132").Value = 11893.6362
$ws.Range("M132").Value = -5254.4645
$ws.Range("N132").Value = -16953.6362
$ws.Range("H136").Value = 58824896
$ws.Range("I136").Value = 83334104
$ws.Range("K136").Value = 250002312
$ws.Range("M136").Value = -249999762

# ---------------------------------------------------------------------------
# Sheet: BSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 50
$ws.Range("K4").Value = 50
$ws.Range("M4").Value = 65
$ws.Range("H82").Value = 4835.087
$ws.Range("I82").Value = 2343.1904
$ws.Range("J82").Value = 31000
$ws.Range("K82").Value = 2343.1904
$ws.Range("L82").Value = 31000
$ws.Range("M82").Value = -1960.1904
$ws.Range("N82").Value = -31766
$ws.Range("H85").Value = 4835.087
$ws.Range("I85").Value = 2343.1904
$ws.Range("J85").Value = 31000
$ws.Range("K85").Value = 2343.1904
$ws.Range("L85").Value = 31000
$ws.Range("M85").Value = -1017.1904
$ws.Range("N85").Value = -33652
$ws.Range("H134").Value = 7832
$ws.Range("I134").Value = 1167.625
$ws.Range("K134").Value = 3502.875
$ws.Range("M134").Value = -967.875

# ---------------------------------------------------------------------------
# Sheet: CRP
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 1009
$ws.Range("I23").Value = 1009
$ws.Range("K23").Value = 1009
$ws.Range("M23").Value = -769
$ws.Range("H27").Value = 1009
$ws.Range("I27").Value = 1009
$ws.Range("K27").Value = 1009
$ws.Range("M27").Value = -817
$ws.Range("H99").Value = 1420.0588
$ws.Range("I99").Value = 1304.6666
$ws.Range("J99").Value = 1549.875
$ws.Range("K99").Value = 1304.6666
$ws.Range("L99").Value = 1549.875
$ws.Range("M99").Value = 193.3334
$ws.Range("N99").Value = -4545.875
$ws.Range("H105").Value = 1075
$ws.Range("I105").Value = 1250
$ws.Range("J105").Value = 900
$ws.Range("K105").Value = 1250
$ws.Range("L105").Value = 900
$ws.Range("M105").Value = 497
$ws.Range("N105").Value = -4394
$ws.Range("H122").Value = 1879.8
$ws.Range("I122").Value = 1499.6666
$ws.Range("J122").Value = 2450
$ws.Range("K122").Value = 4498.9998
$ws.Range("L122").Value = 7350
$ws.Range("M122").Value = -2048.9998
$ws.Range("N122").Value = -12250
$ws.Range("H126").Value = 1420.0588
$ws.Range("I126").Value = 1304.6666
$ws.Range("J126").Value = 1549.875
$ws.Range("K126").Value = 3913.9998
$ws.Range("L126").Value = 4649.625
$ws.Range("M126").Value = -1443.9998
$ws.Range("N126").Value = -9589.625

# ---------------------------------------------------------------------------
# Sheet: CUL
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 746.6
$ws.Range("I97").Value = 800
$ws.Range("J97").Value = 711
$ws.Range("K97").Value = 2400
$ws.Range("L97").Value = 2133
$ws.Range("M97").Value = -1904
$ws.Range("N97").Value = -3125

# ---------------------------------------------------------------------------
# Sheet: GSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4050
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 4050
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 4050
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -6046
$ws.Range("H83").Value = 4050
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 4050
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 20250
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -30234

# ---------------------------------------------------------------------------
# Sheet: LTW
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1986.6666
$ws.Range("I7").Value = 1625.7142
$ws.Range("J7").Value = 3250
$ws.Range("K7").Value = 1625.7142
$ws.Range("L7").Value = 3250
$ws.Range("M7").Value = -1513.7142
$ws.Range("N7").Value = -3474
$ws.Range("H40").Value = 2332.1875
$ws.Range("I40").Value = 2143.8333
$ws.Range("J40").Value = 2897.25
$ws.Range("K40").Value = 2143.8333
$ws.Range("L40").Value = 2897.25
$ws.Range("M40").Value = -2007.8333
$ws.Range("N40").Value = -3169.25
$ws.Range("H82").Value = 1132.2222
$ws.Range("J82").Value = 1057.1428
$ws.Range("L82").Value = 1057.1428
$ws.Range("N82").Value = -1779.1428
$ws.Range("H85").Value = 1132.2222
$ws.Range("J85").Value = 1057.1428
$ws.Range("L85").Value = 1057.1428
$ws.Range("N85").Value = -3553.1428
$ws.Range("H122").Value = 19233316
$ws.Range("I122").Value = 20835758
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 62507274
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -62504824
$ws.Range("N122").Value = -16900
$ws.Range("H126").Value = 1986.6666
$ws.Range("I126").Value = 1625.7142
$ws.Range("J126").Value = 3250
$ws.Range("K126").Value = 4877.142599999999
$ws.Range("L126").Value = 9750
$ws.Range("M126").Value = -2407.142599999999
$ws.Range("N126").Value = -14690
$ws.Range("H130").Value = 34664.5
$ws.Range("J130").Value = 34664.5
$ws.Range("L130").Value = 34664.5
$ws.Range("N130").Value = -44704.5

# ---------------------------------------------------------------------------
# Sheet: WVR
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 12502547
$ws.Range("I122").Value = 14708456
$ws.Range("J122").Value = 2400
$ws.Range("K122").Value = 44125368
$ws.Range("L122").Value = 7200
$ws.Range("M122").Value = -44122918
$ws.Range("N122").Value = -12100
$ws.Range("H126").Value = 142858560
$ws.Range("J126").Value = 880
$ws.Range("L126").Value = 2640
$ws.Range("N126").Value = -7580
